$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the version and "last tested with" text, per the latest release
$ws.Range("A5").Value2 = "(Version: 1.0.1)"
$ws.Range("A6").Value2 = "(Last tested with: ReportServer 4.0.0-6053) "

# Move the active selection from A5 to A6
$ws.Range("A6").Select()
